$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the "ECs -> Resolving-Mac" edge row no longer present)
$ws.Rows.Item(3).Delete()

# Update row 2 with the new TPM-derived values, now the only data row
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 66.00836066666666
$ws.Range("H2").Value = 198.025082
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5468886666666667
$ws.Range("N2").Value = 1.640666
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 36.09922435384578
$ws.Range("R2").Value = 324.893019184612
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
